$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing error-message cells to "Done" (rows 35, 38, 39, 40) ---
$ws.Range("E35").Value = "Done"
$ws.Range("E38").Value = "Done"
$ws.Range("E39").Value = "Done"
$ws.Range("E40").Value = "Done"

# --- Fill previously-empty Status cells with "Done" (rows 41-49) ---
for ($r = 41; $r -le 49; $r++) {
    $ws.Cells.Item($r, 5).Value = "Done"
}

# --- Append new rows 50-72 ---
# Use row 49 (A:F) as a formatting/shape template for the blank B/C/F cells,
# then overwrite A/D/E with the real data for each new row.
$newRowsData = @(
    @{Row=50; Url="https://www.youtube.com/watch?v=70y6hMyRZfQ&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=25&pp=iAQB"; Status="Done"},
    @{Row=51; Url="https://www.youtube.com/watch?v=eHiqrRloaGc&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=28&pp=iAQB"; Status="Done"},
    @{Row=52; Url="https://www.youtube.com/watch?v=Kp61Z_TKJFk&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=31&pp=iAQB"; Status=""},
    @{Row=53; Url="https://www.youtube.com/watch?v=p6WfKfrLhdQ&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=34&pp=iAQB"; Status=""},
    @{Row=54; Url="https://www.youtube.com/watch?v=pswBnrF5Z5Y&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=37&pp=iAQB"; Status=""},
    @{Row=55; Url="https://www.youtube.com/watch?v=PUscmv4YqMU&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=42&pp=iAQB"; Status=""},
    @{Row=56; Url="https://www.youtube.com/watch?v=viY5RFtIjo0&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=46&pp=iAQB"; Status=""},
    @{Row=57; Url="https://www.youtube.com/watch?v=Z4EpFOmyf_g&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=49&pp=iAQB"; Status=""},
    @{Row=58; Url="https://www.youtube.com/watch?v=K1ipCo_KCdI&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=55&pp=iAQB"; Status=""},
    @{Row=59; Url="https://www.youtube.com/watch?v=gunJ_nYCe4k&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=59&pp=iAQB"; Status=""},
    @{Row=60; Url="https://www.youtube.com/watch?v=ZYfAO1VJ0r8&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=63&pp=iAQB"; Status=""},
    @{Row=61; Url="https://www.youtube.com/watch?v=5-Hi6SXvrgU&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=67&pp=iAQB"; Status=""},
    @{Row=62; Url="https://www.youtube.com/watch?v=NALYks6kRM0&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=71&pp=iAQB"; Status=""},
    @{Row=63; Url="https://www.youtube.com/watch?v=AFvQjxICn3E&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=75&pp=iAQB"; Status=""},
    @{Row=64; Url="https://www.youtube.com/watch?v=aNPW8VA8wBo&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=79&pp=iAQB"; Status=""},
    @{Row=65; Url="https://www.youtube.com/watch?v=F5jyrHQff9w&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=82&pp=iAQB"; Status=""},
    @{Row=66; Url="https://www.youtube.com/watch?v=hKVpdmjvGRA&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=86&pp=iAQB"; Status=""},
    @{Row=67; Url="https://www.youtube.com/watch?v=8OLHwss-D20&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=90&pp=iAQB"; Status=""},
    @{Row=68; Url="https://www.youtube.com/watch?v=eSs9wOgLTWE&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=94&pp=iAQB"; Status=""},
    @{Row=69; Url="https://www.youtube.com/watch?v=7d5Vs-ha1XU&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=97&pp=iAQB"; Status=""},
    @{Row=70; Url="https://www.youtube.com/watch?v=QFkk5cN91B0&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=100&pp=iAQB"; Status=""},
    @{Row=71; Url="https://www.youtube.com/watch?v=gQQG8KK5gQU&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=105&pp=iAQB"; Status=""},
    @{Row=72; Url="https://www.youtube.com/watch?v=2uj0bE7yQhw&list=PLipLTGfSGq9EZDjrAi7VnUf7loK35l19x&index=109&pp=iAQB"; Status=""}
)

foreach ($item in $newRowsData) {
    $r = $item.Row
    # Seed the row with the same blank shape as row 49 (keeps B/C/F blank cells)
    $ws.Range("A49:F49").Copy($ws.Range("A" + $r + ":F" + $r))

    $ws.Cells.Item($r, 1).Value = $item.Url
    $ws.Cells.Item($r, 4).Value = 1
    if ($item.Status -ne "") {
        $ws.Cells.Item($r, 5).Value = $item.Status
    } else {
        $ws.Cells.Item($r, 5).ClearContents()
    }
}
